$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 27: the "note" cell (D27) gets a new value -------------------------
$ws.Range("D27").Value = "in private 0,00473"

# --- Row 28 (new) ------------------------------------------------------------
$ws.Range("A28").Value = "hybrid_cbf_cfAdjCosine_w0.13cf_w0.87cbf_popularity500"
# B28 looks like a plain number ("0.00588") so force text first, otherwise
# Excel would silently store it as a Number instead of matching the sheet's
# existing "numeric text" convention for this column.
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "0.00588"
$ws.Range("B28").NumberFormat = "0.00000"
$ws.Range("B28").HorizontalAlignment = -4131
$ws.Range("C28").Value = "…"
$ws.Range("D28").Value = "in private 0,00601"

# --- Row 29 (new) ------------------------------------------------------------
$ws.Range("A29").Value = "hybrid_cbf_cf_w0.13cf_w0.87cbf_popularity1000"
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "0.00573"
$ws.Range("B29").NumberFormat = "0.00000"
$ws.Range("B29").HorizontalAlignment = -4131
$ws.Range("C29").Value = "…"
$ws.Range("D29").Value = "in private 0,00548"

# --- Match formatting used elsewhere in the sheet for these columns --------
$ws.Range("D27").Interior.Color = 15773696
$ws.Range("C28:C29").NumberFormat = "0.00000"
$ws.Range("C28:C29").HorizontalAlignment = -4131
$ws.Range("D28:D29").NumberFormat = "0.00"
$ws.Range("D28:D29").Interior.Color = 15773696
$ws.Range("A28:A29").VerticalAlignment = -4108

# --- Update selection / scroll position to match the new view state --------
[void]$ws.Range("D26").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 2
